{"js": "// Remove the trailing \"Ver no Jupiter...\" and \"\u00a9 2020 ...\" paragraphs\n// from the end of the document body (site-footer boilerplate stripped\n// from the page at build time), leaving the rest of the document intact.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst targets = [\n  \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n  \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n];\n\nfor (const paragraph of paragraphs.items) {\n  if (targets.includes(paragraph.text.trim())) {\n    paragraph.delete();\n  }\n}\n\nawait context.sync();\n", "ps1": "# Remove the trailing \"Ver no Jupiter...\" and \"(c) 2020 ...\" site-footer\n# paragraphs from the end of the document, leaving the rest untouched.\n$d = $word.ActiveDocument\n\n$targets = @(\n    \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n    [char]0x00A9 + \" 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n)\n\nfor ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {\n    $para = $d.Paragraphs.Item($i)\n    $text = $para.Range.Text.TrimEnd(\"`r\", \"`a\")\n    if ($targets -contains $text) {\n        $para.Range.Delete()\n    }\n}\n"}
